$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Handback report generation: for each language sheet (zh-cn, de-de), the
# first two data rows (a.md, b.md) move from "Ready for handoff" to
# "Handed back: in sync with en-US". Two new columns of data are populated
# for those rows - "Latest Target File" (E) and "Latest Handback File" (F) -
# and "Latest Handback DateTime" (G) is stamped with the handback time.
# ---------------------------------------------------------------------------

function Set-HandbackRow {
    param($ws, $row, $targetDisplay, $targetUrl, $handbackDisplay, $handbackUrl, $handbackDateTime)

    # Status -> Handed back
    $ws.Cells.Item($row, 2).Value = "Handed back: in sync with en-US"

    # E: Latest Target File (new hyperlink cell)
    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value = $targetDisplay
    $ws.Hyperlinks.Add($eCell, $targetUrl, [Type]::Missing, [Type]::Missing, $targetDisplay) | Out-Null
    $eCell.Font.Underline = 2
    $eCell.Font.Color = 15570276

    # F: Latest Handback File (new hyperlink cell)
    $fCell = $ws.Cells.Item($row, 6)
    $fCell.Value = $handbackDisplay
    $ws.Hyperlinks.Add($fCell, $handbackUrl, [Type]::Missing, [Type]::Missing, $handbackDisplay) | Out-Null
    $fCell.Font.Underline = 2
    $fCell.Font.Color = 15570276

    # G: Latest Handback DateTime
    $ws.Cells.Item($row, 7).Value = $handbackDateTime
}

# --------------------------- zh-cn sheet -----------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

Set-HandbackRow `
    $wsZh `
    2 `
    "a.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/ba8b4f10e8cb154333ecc169904bb55277c0f3ba/e2e/a.md" `
    "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c4813d29c9186abbb137299024e8129fbb20939c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" `
    "2016-02-22 09:16:26"

Set-HandbackRow `
    $wsZh `
    3 `
    "a.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/ba8b4f10e8cb154333ecc169904bb55277c0f3ba/e2e/a.md" `
    "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c4813d29c9186abbb137299024e8129fbb20939c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" `
    "2016-02-22 09:16:26"

# --------------------------- de-de sheet -----------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

Set-HandbackRow `
    $wsDe `
    2 `
    "a.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/ba8b4f10e8cb154333ecc169904bb55277c0f3ba/e2e/a.md" `
    "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bdb75358cedcaa59bff617c079c464a5e61dc89f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" `
    "2016-02-22 09:16:49"

Set-HandbackRow `
    $wsDe `
    3 `
    "a.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/ba8b4f10e8cb154333ecc169904bb55277c0f3ba/e2e/a.md" `
    "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bdb75358cedcaa59bff617c079c464a5e61dc89f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" `
    "2016-02-22 09:16:49"
